# Insert a new row above the current row 3 ("We are doing our best...") that
# explains orders must be state-wide, pushing the existing rows 3-8 down to
# rows 4-9. The new cell keeps the worksheet's default (unstyled) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3").Insert()
$ws.Range("A3").Value = "We only record orders if they apply state-wide"
$ws.Range("A3").ClearFormats()

$ws.Range("A4").Select()
